{"js": "// Replace the 25 two-digit multiplication problems in the worksheet table\n// with a new set of problems, per the commit's regenerated answer key.\nconst replacements = [\n  [\"17\u00d742=714\", \"94\u00d752=4888\"],\n  [\"57\u00d793=5301\", \"29\u00d776=2204\"],\n  [\"49\u00d761=2989\", \"37\u00d779=2923\"],\n  [\"47\u00d728=1316\", \"36\u00d740=1440\"],\n  [\"63\u00d791=5733\", \"60\u00d728=1680\"],\n  [\"66\u00d793=6138\", \"45\u00d763=2835\"],\n  [\"12\u00d714=168\", \"49\u00d729=1421\"],\n  [\"34\u00d752=1768\", \"37\u00d770=2590\"],\n  [\"50\u00d779=3950\", \"69\u00d787=6003\"],\n  [\"21\u00d732=672\", \"52\u00d747=2444\"],\n  [\"92\u00d755=5060\", \"60\u00d784=5040\"],\n  [\"40\u00d743=1720\", \"28\u00d771=1988\"],\n  [\"30\u00d717=510\", \"97\u00d727=2619\"],\n  [\"66\u00d779=5214\", \"67\u00d794=6298\"],\n  [\"95\u00d752=4940\", \"19\u00d717=323\"],\n  [\"94\u00d785=7990\", \"24\u00d720=480\"],\n  [\"69\u00d711=759\", \"11\u00d738=418\"],\n  [\"62\u00d776=4712\", \"77\u00d717=1309\"],\n  [\"67\u00d795=6365\", \"78\u00d738=2964\"],\n  [\"83\u00d750=4150\", \"26\u00d738=988\"],\n  [\"48\u00d796=4608\", \"75\u00d782=6150\"],\n  [\"15\u00d716=240\", \"87\u00d744=3828\"],\n  [\"26\u00d778=2028\", \"99\u00d731=3069\"],\n  [\"43\u00d778=3354\", \"12\u00d764=768\"],\n  [\"84\u00d746=3864\", \"30\u00d786=2580\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(`Expected exactly 1 match for \"${oldText}\" but found ${results.items.length}`);\n  }\n  results.items[0].insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 two-digit multiplication problems in the worksheet table\n# with a new set of problems, per the commit's regenerated answer key.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"17\u00d742=714\", \"94\u00d752=4888\"),\n  @(\"57\u00d793=5301\", \"29\u00d776=2204\"),\n  @(\"49\u00d761=2989\", \"37\u00d779=2923\"),\n  @(\"47\u00d728=1316\", \"36\u00d740=1440\"),\n  @(\"63\u00d791=5733\", \"60\u00d728=1680\"),\n  @(\"66\u00d793=6138\", \"45\u00d763=2835\"),\n  @(\"12\u00d714=168\", \"49\u00d729=1421\"),\n  @(\"34\u00d752=1768\", \"37\u00d770=2590\"),\n  @(\"50\u00d779=3950\", \"69\u00d787=6003\"),\n  @(\"21\u00d732=672\", \"52\u00d747=2444\"),\n  @(\"92\u00d755=5060\", \"60\u00d784=5040\"),\n  @(\"40\u00d743=1720\", \"28\u00d771=1988\"),\n  @(\"30\u00d717=510\", \"97\u00d727=2619\"),\n  @(\"66\u00d779=5214\", \"67\u00d794=6298\"),\n  @(\"95\u00d752=4940\", \"19\u00d717=323\"),\n  @(\"94\u00d785=7990\", \"24\u00d720=480\"),\n  @(\"69\u00d711=759\", \"11\u00d738=418\"),\n  @(\"62\u00d776=4712\", \"77\u00d717=1309\"),\n  @(\"67\u00d795=6365\", \"78\u00d738=2964\"),\n  @(\"83\u00d750=4150\", \"26\u00d738=988\"),\n  @(\"48\u00d796=4608\", \"75\u00d782=6150\"),\n  @(\"15\u00d716=240\", \"87\u00d744=3828\"),\n  @(\"26\u00d778=2028\", \"99\u00d731=3069\"),\n  @(\"43\u00d778=3354\", \"12\u00d764=768\"),\n  @(\"84\u00d746=3864\", \"30\u00d786=2580\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $rng = $d.Content\n  $rng.Find.ClearFormatting()\n  $rng.Find.Execute($oldText, $true, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
